$wb = $excel.ActiveWorkbook

# Sheet 1: Spherical
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = [double]"0.7692941849576688"
$ws.Range("B2").Value = [double]"1.542730847490668"
$ws.Range("C2").Value = [double]"1.490879221317486"
$ws.Range("D2").Value = [double]"1.242067167060891"
$ws.Range("F2").Value = [double]"2.075685827223915e-29"
$ws.Range("G2").Value = [double]"1.818878597560449e-14"
$ws.Range("H2").Value = [double]"4.555969520556426e-15"
$ws.Range("I2").Value = [double]"0.7739693040866447"
$ws.Range("J2").Value = [double]"0.8084986932692988"
$ws.Range("K2").Value = [double]"3.364388869890441"
$ws.Range("L2").Value = [double]"0.8991655538716431"
$ws.Range("M2").Value = [double]"0.7819351888578722"
$ws.Range("N2").Value = [double]"1.353359040557771"
$ws.Range("O2").Value = [double]"0.9574498701462588"
$ws.Range("P2").Value = [double]"1.163339606717562"

# Sheet 2: Gaussian
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = [double]"0.8276108165723453"
$ws.Range("B2").Value = [double]"1.152767263359928"
$ws.Range("C2").Value = [double]"1.93271542165783"
$ws.Range("D2").Value = [double]"1.073669997420031"
$ws.Range("F2").Value = [double]"2.549350771083462e-29"
$ws.Range("G2").Value = [double]"1.9067849121359e-14"
$ws.Range("H2").Value = [double]"5.049109595843074e-15"
$ws.Range("I2").Value = [double]"0.812215503691748"
$ws.Range("J2").Value = [double]"0.6716942549239149"
$ws.Range("K2").Value = [double]"0.4175751219327417"
$ws.Range("L2").Value = [double]"0.8195695546589776"
$ws.Range("M2").Value = [double]"0.837852471077039"
$ws.Range("N2").Value = [double]"1.006323867765007"
$ws.Range("O2").Value = [double]"1.017921705836737"
$ws.Range("P2").Value = [double]"1.0031569507136"

# Sheet 3: Exponential
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = [double]"0.7435723829541865"
$ws.Range("B2").Value = [double]"1.714732655925962"
$ws.Range("C2").Value = [double]"1.319891597474623"
$ws.Range("D2").Value = [double]"1.30947800895088"
$ws.Range("F2").Value = [double]"1.705243987628719e-29"
$ws.Range("G2").Value = [double]"1.152722558947848e-14"
$ws.Range("H2").Value = [double]"4.129459998145906e-15"
$ws.Range("I2").Value = [double]"0.7613785966610122"
$ws.Range("J2").Value = [double]"0.8535349236796232"
$ws.Range("K2").Value = [double]"3.483026505441671"
$ws.Range("L2").Value = [double]"0.9238695382355797"
$ws.Range("M2").Value = [double]"0.7554879657324924"
$ws.Range("N2").Value = [double]"1.517496428552264"
$ws.Range("O2").Value = [double]"0.796816168213166"
$ws.Range("P2").Value = [double]"1.23186704986872"

# Sheet 4: Linear
$ws = $wb.Worksheets.Item(4)
$ws.Range("A2").Value = [double]"0.2813534431453613"
$ws.Range("B2").Value = [double]"4.805592834750094"
$ws.Range("C2").Value = [double]"1.732046185857099"
$ws.Range("D2").Value = [double]"2.192166242498523"
$ws.Range("F2").Value = [double]"2.009116636475154e-30"
$ws.Range("G2").Value = [double]"3.955089848909272e-15"
$ws.Range("H2").Value = [double]"1.417433115344478e-15"
$ws.Range("I2").Value = [double]"0.3845846849934301"
$ws.Range("J2").Value = [double]"2.201304897948272"
$ws.Range("K2").Value = [double]"4.455349818907312"
$ws.Range("L2").Value = [double]"1.483679513219843"
$ws.Range("M2").Value = [double]"0.2974918202433614"
$ws.Range("N2").Value = [double]"4.359923048381081"
$ws.Range("O2").Value = [double]"2.664742142868901"
$ws.Range("P2").Value = [double]"2.088042875129982"

# Sheet 5: Power
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = [double]"0.5477659675530802"
$ws.Range("B2").Value = [double]"3.0240910573188"
$ws.Range("C2").Value = [double]"1.687083469284017"
$ws.Range("D2").Value = [double]"1.738991390812157"
$ws.Range("F2").Value = [double]"3.874062009173368e-30"
$ws.Range("G2").Value = [double]"5.612399192364684e-15"
$ws.Range("H2").Value = [double]"1.968263704175172e-15"
$ws.Range("I2").Value = [double]"0.613662684371109"
$ws.Range("J2").Value = [double]"1.381906176879894"
$ws.Range("K2").Value = [double]"4.770316126117107"
$ws.Range("L2").Value = [double]"1.175545055231782"
$ws.Range("M2").Value = [double]"0.5543866075244759"
$ws.Range("N2").Value = [double]"2.765576482247305"
$ws.Range("O2").Value = [double]"1.814180411310868"
$ws.Range("P2").Value = [double]"1.663002249621841"

# Sheet 6: HoleEffect
$ws = $wb.Worksheets.Item(6)
$ws.Range("A2").Value = [double]"-0.005407713470203523"
$ws.Range("B2").Value = [double]"6.723166009452089"
$ws.Range("C2").Value = [double]"2.063041572240144"
$ws.Range("D2").Value = [double]"2.592906864785561"
$ws.Range("F2").Value = [double]"6.423361550520309e-30"
$ws.Range("G2").Value = [double]"8.227612281004692e-15"
$ws.Range("H2").Value = [double]"2.534435154135988e-15"
$ws.Range("I2").Value = [double]"0.01501071341073656"
$ws.Range("J2").Value = [double]"3.523249565169453"
$ws.Range("K2").Value = [double]"3.212357516838896"
$ws.Range("M2").Value = [double]"-0.001948746461204243"
$ws.Range("N2").Value = [double]"6.218318247206794"
$ws.Range("O2").Value = [double]"2.590983844690503"
$ws.Range("P2").Value = [double]"2.493655599156947"
